# Update the Compliance Waste Returns title strings to include a space
# before the {0} placeholder, and update the saved cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B25 = DMSExcelReturnTitle value: "Compliance Waste Returns{0} {1}"
#       -> "Compliance Waste Returns {0} {1}"
$ws.Range("B25").Value = "Compliance Waste Returns {0} {1}"

# B23 = DMSEmailTitle value: "Compliance Waste Returns{0} {1} - Email and Submission"
#       -> "Compliance Waste Returns {0} {1} - Email and Submission"
$ws.Range("B23").Value = "Compliance Waste Returns {0} {1} - Email and Submission"

# Update the sheet's active selection from B22 to B26
$ws.Range("B26").Select()
